# SET DATA FROM ARRAY INTO CELLS
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The workbook used to carry a duplicate "Copy Sheet" tab - drop it so only
# the original "First Sheet" remains.
$wb.Worksheets.Item("Copy Sheet").Delete()

$ws = $wb.Worksheets.Item("First Sheet")
$ws.Activate()

# Remove the old "Hello World!" text that lived in A1.
$ws.Range("A1").ClearContents()

# Build the array of numbers and drop it straight into column B (B1:B3).
$values = @(100, 53, 86)
$data = New-Object 'object[,]' $values.Length,1
for ($i = 0; $i -lt $values.Length; $i++) {
    $data[$i, 0] = $values[$i]
}
$ws.Range("B1:B3").Value = $data
